$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.567170977592468
$ws.Range("B1").Value = 3.672008514404297
$ws.Range("C1").Value = 5.553871154785156
$ws.Range("D1").Value = 1.376509070396423
$ws.Range("E1").Value = 0.8022044897079468
